$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1938.75
$ws.Range("I28").Value = 2285
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 2285
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -1800
$ws.Range("N28").Value = -1870

$ws.Range("H40").Value = 3149.9

$ws.Range("H70").Value = 10109005
$ws.Range("I70").Value = 22223850
$ws.Range("K70").Value = 66671550
$ws.Range("M70").Value = -66671280

$ws.Range("H73").Value = 10109005
$ws.Range("I73").Value = 22223850
$ws.Range("K73").Value = 66671550
$ws.Range("M73").Value = -66670614

$ws.Range("H80").Value = 1406935.8
$ws.Range("J80").Value = 7238.5
$ws.Range("L80").Value = 21715.5
$ws.Range("N80").Value = -23711.5

$ws.Range("H83").Value = 1406935.8
$ws.Range("J83").Value = 7238.5
$ws.Range("L83").Value = 65146.5
$ws.Range("N83").Value = -75130.5

$ws.Range("H100").Value = 8407.5
$ws.Range("I100").Value = 9722.5
$ws.Range("K100").Value = 9722.5
$ws.Range("M100").Value = -9181.5

$ws.Range("H103").Value = 1275.6666
$ws.Range("I103").Value = 1787
$ws.Range("K103").Value = 5361
$ws.Range("M103").Value = -4775

$ws.Range("H106").Value = 2221.7856
$ws.Range("I106").Value = 2476.0833
$ws.Range("K106").Value = 2476.0833
$ws.Range("M106").Value = -1845.0833

$ws.Range("H111").Value = 1593.75
$ws.Range("I111").Value = 1463
$ws.Range("J111").Value = 1753.5555
$ws.Range("K111").Value = 4389
$ws.Range("L111").Value = 5260.666499999999
$ws.Range("M111").Value = -1322
$ws.Range("N111").Value = -11394.6665

$ws.Range("H113").Value = 4512.1055
$ws.Range("I113").Value = 3544.8
$ws.Range("J113").Value = 5586.8887
$ws.Range("K113").Value = 3544.8
$ws.Range("L113").Value = 5586.8887
$ws.Range("M113").Value = -290.8000000000002
$ws.Range("N113").Value = -12094.8887

$ws.Range("H132").Value = 4951517
$ws.Range("I132").Value = 4951517
$ws.Range("K132").Value = 14854551
$ws.Range("M132").Value = -14852021

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28509.65
$ws.Range("I32").Value = 33141.324
$ws.Range("K32").Value = 33141.324
$ws.Range("M32").Value = -32854.324

$ws.Range("H45").Value = 3345.25
$ws.Range("I45").Value = 976.875
$ws.Range("K45").Value = 976.875
$ws.Range("M45").Value = -599.875

$ws.Range("H102").Value = 864.13336
$ws.Range("I102").Value = 680.25
$ws.Range("K102").Value = 680.25
$ws.Range("M102").Value = 941.75

$ws.Range("H132").Value = 1660.1428
$ws.Range("I132").Value = 1363.2667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4089.800099999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1559.800099999999
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4365.625
$ws.Range("I16").Value = 4154.3335
$ws.Range("J16").Value = 4999.5
$ws.Range("K16").Value = 4154.3335
$ws.Range("L16").Value = 4999.5
$ws.Range("M16").Value = -3867.3335
$ws.Range("N16").Value = -5573.5

$ws.Range("H62").Value = 5648.8
$ws.Range("J62").Value = 5500
$ws.Range("L62").Value = 5500
$ws.Range("N62").Value = -6748

$ws.Range("H65").Value = 5648.8
$ws.Range("J65").Value = 5500
$ws.Range("L65").Value = 27500
$ws.Range("N65").Value = -33740

$ws.Range("H94").Value = 1051.72
$ws.Range("I94").Value = 1255.8889
$ws.Range("J94").Value = 936.875
$ws.Range("K94").Value = 1255.8889
$ws.Range("L94").Value = 936.875
$ws.Range("M94").Value = -804.8888999999999
$ws.Range("N94").Value = -1838.875

$ws.Range("H105").Value = 1745.5
$ws.Range("I105").Value = 977.25
$ws.Range("K105").Value = 977.25
$ws.Range("M105").Value = 769.75

$ws.Range("H113").Value = 4365.625
$ws.Range("I113").Value = 4154.3335
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 4154.3335
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -1984.3335
$ws.Range("N113").Value = -9339.5

$ws.Range("H132").Value = 43743.625
$ws.Range("I132").Value = 51517.9
$ws.Range("K132").Value = 154553.7
$ws.Range("M132").Value = -152023.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3949225.2
$ws.Range("J11").Value = 999.5
$ws.Range("L11").Value = 2998.5
$ws.Range("N11").Value = -3278.5

$ws.Range("H131").Value = 2278.2307
$ws.Range("I131").Value = 3130.8
$ws.Range("J131").Value = 1932.5946
$ws.Range("K131").Value = 9392.400000000001
$ws.Range("L131").Value = 5797.783799999999
$ws.Range("M131").Value = -4352.400000000001
$ws.Range("N131").Value = -15877.7838

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2631.3635
$ws.Range("I80").Value = 2775.25
$ws.Range("K80").Value = 2775.25
$ws.Range("M80").Value = -1777.25

$ws.Range("H83").Value = 2631.3635
$ws.Range("I83").Value = 2775.25
$ws.Range("K83").Value = 13876.25
$ws.Range("M83").Value = -8884.25

$ws.Range("H113").Value = 1318.2632
$ws.Range("I113").Value = 877.36365
$ws.Range("K113").Value = 877.36365
$ws.Range("M113").Value = 1292.63635

$ws.Range("I122").Value = 2060.6365
$ws.Range("K122").Value = 6181.9095
$ws.Range("M122").Value = -3731.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3231.4666
$ws.Range("I40").Value = 3041.0833
$ws.Range("J40").Value = 3993
$ws.Range("K40").Value = 3041.0833
$ws.Range("L40").Value = 3993
$ws.Range("M40").Value = -2905.0833
$ws.Range("N40").Value = -4265

$ws.Range("H61").Value = 1675.0435
$ws.Range("I61").Value = 1707.8422
$ws.Range("J61").Value = 1519.25
$ws.Range("K61").Value = 1707.8422
$ws.Range("L61").Value = 1519.25
$ws.Range("M61").Value = -1505.8422
$ws.Range("N61").Value = -1923.25

$ws.Range("H82").Value = 2959.6785
$ws.Range("I82").Value = 2769.611
$ws.Range("K82").Value = 2769.611
$ws.Range("M82").Value = -2408.611

$ws.Range("H85").Value = 2959.6785
$ws.Range("I85").Value = 2769.611
$ws.Range("K85").Value = 2769.611
$ws.Range("M85").Value = -1521.611

$ws.Range("H100").Value = 2706.7
$ws.Range("I100").Value = 2563
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2563
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2022
$ws.Range("N100").Value = -5082

$ws.Range("H113").Value = 1675.0435
$ws.Range("I113").Value = 1707.8422
$ws.Range("J113").Value = 1519.25
$ws.Range("K113").Value = 1707.8422
$ws.Range("L113").Value = 1519.25
$ws.Range("M113").Value = 462.1578
$ws.Range("N113").Value = -5859.25

$ws.Range("H122").Value = 1691.4166
$ws.Range("I122").Value = 1710.2222
$ws.Range("J122").Value = 1635
$ws.Range("K122").Value = 5130.6666
$ws.Range("L122").Value = 4905
$ws.Range("M122").Value = -2680.6666
$ws.Range("N122").Value = -9805

$ws.Range("H132").Value = 4829.0347
$ws.Range("I132").Value = 4125.5557
$ws.Range("K132").Value = 12376.6671
$ws.Range("M132").Value = -9846.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 872.72
$ws.Range("I113").Value = 901
$ws.Range("J113").Value = 547.5
$ws.Range("K113").Value = 2703
$ws.Range("L113").Value = 1642.5
$ws.Range("M113").Value = -533
$ws.Range("N113").Value = -5982.5

$ws.Range("H136").Value = 15943.547
$ws.Range("I136").Value = 20405.703
$ws.Range("J136").Value = 5624.8125
$ws.Range("K136").Value = 61217.109
$ws.Range("L136").Value = 16874.4375
$ws.Range("M136").Value = -58667.109
$ws.Range("N136").Value = -21974.4375
